$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cD = $ws.Cells.Item(2, 4)
$cD.NumberFormat = "@"
$cD.Value = "27.184.36"
$cD.Style = "Normal"
$ws.Cells.Item(2, 5).Value = "  +1.00%  "
$cD = $ws.Cells.Item(3, 4)
$cD.NumberFormat = "@"
$cD.Value = "1.857.84"
$cD.Style = "Normal"
$ws.Cells.Item(3, 5).Value = "  +0.62%  "
$ws.Cells.Item(4, 5).Value = "  +1.13%  "
$ws.Cells.Item(5, 5).Value = "  +1.10%  "
$cD = $ws.Cells.Item(6, 4)
$cD.NumberFormat = "@"
$cD.Value = "311.28"
$cD.Style = "Normal"
$ws.Cells.Item(6, 5).Value = "  +0.64%  "
$cD = $ws.Cells.Item(7, 4)
$cD.NumberFormat = "@"
$cD.Value = "0.4794"
$cD.Style = "Normal"
$ws.Cells.Item(7, 5).Value = "  +2.03%  "
$cD = $ws.Cells.Item(8, 4)
$cD.NumberFormat = "@"
$cD.Value = "0.3713"
$cD.Style = "Normal"
$ws.Cells.Item(8, 5).Value = "  +1.34%  "
$cD = $ws.Cells.Item(9, 4)
$cD.NumberFormat = "@"
$cD.Value = "0.07288"
$cD.Style = "Normal"
$ws.Cells.Item(9, 5).Value = "  +1.77%  "
$cD = $ws.Cells.Item(10, 4)
$cD.NumberFormat = "@"
$cD.Value = "0.9343"
$cD.Style = "Normal"
$ws.Cells.Item(10, 5).Value = "  +0.82%  "
$cD = $ws.Cells.Item(11, 4)
$cD.NumberFormat = "@"
$cD.Value = "20.08"
$cD.Style = "Normal"
$ws.Cells.Item(11, 5).Value = "  +2.49%  "
$cD = $ws.Cells.Item(12, 4)
$cD.NumberFormat = "@"
$cD.Value = "0.07861"
$cD.Style = "Normal"
$ws.Cells.Item(12, 5).Value = "  +2.00%  "
$cD = $ws.Cells.Item(13, 4)
$cD.NumberFormat = "@"
$cD.Value = "1.863.27"
$cD.Style = "Normal"
$ws.Cells.Item(13, 5).Value = "  +1.16%  "
$cD = $ws.Cells.Item(14, 4)
$cD.NumberFormat = "@"
$cD.Value = "5.419"
$cD.Style = "Normal"
$ws.Cells.Item(14, 5).Value = "  +2.57%  "
$cD = $ws.Cells.Item(15, 4)
$cD.NumberFormat = "@"
$cD.Value = "6.536"
$cD.Style = "Normal"
$ws.Cells.Item(15, 5).Value = "  +1.98%  "
$cD = $ws.Cells.Item(16, 4)
$cD.NumberFormat = "@"
$cD.Value = "89.98"
$cD.Style = "Normal"
$ws.Cells.Item(16, 5).Value = "  +1.86%  "
$cD = $ws.Cells.Item(17, 4)
$cD.NumberFormat = "@"
$cD.Value = "1.020"
$cD.Style = "Normal"
$ws.Cells.Item(17, 5).Value = "  +1.18%  "
$cD = $ws.Cells.Item(18, 4)
$cD.NumberFormat = "@"
$cD.Value = "0.000008722"
$cD.Style = "Normal"
$ws.Cells.Item(18, 5).Value = "  +0.90%  "
$cD = $ws.Cells.Item(19, 4)
$cD.NumberFormat = "@"
$cD.Value = "1.018"
$cD.Style = "Normal"
$ws.Cells.Item(19, 5).Value = "  +1.15%  "
$cD = $ws.Cells.Item(20, 4)
$cD.NumberFormat = "@"
$cD.Value = "27.230.44"
$cD.Style = "Normal"
$ws.Cells.Item(20, 5).Value = "  +1.04%  "
$cD = $ws.Cells.Item(21, 4)
$cD.NumberFormat = "@"
$cD.Value = "14.69"
$cD.Style = "Normal"
$ws.Cells.Item(21, 5).Value = "  +1.70%  "
$ws.Cells.Item(22, 5).Value = "  +1.37%  "
$cD = $ws.Cells.Item(23, 4)
$cD.NumberFormat = "@"
$cD.Value = "10.66"
$cD.Style = "Normal"
$ws.Cells.Item(23, 5).Value = "  +0.53%  "
$cD = $ws.Cells.Item(24, 4)
$cD.NumberFormat = "@"
$cD.Value = "1.951"
$cD.Style = "Normal"
$ws.Cells.Item(24, 5).Value = "  +0.98%  "
$cD = $ws.Cells.Item(25, 4)
$cD.NumberFormat = "@"
$cD.Value = "153.57"
$cD.Style = "Normal"
$ws.Cells.Item(25, 5).Value = "  +1.05%  "
$cD = $ws.Cells.Item(26, 4)
$cD.NumberFormat = "@"
$cD.Value = "18.46"
$cD.Style = "Normal"
$ws.Cells.Item(26, 5).Value = "  +1.16%  "
$ws.Cells.Item(27, 5).Value = "  -0.91%  "
$cD = $ws.Cells.Item(28, 4)
$cD.NumberFormat = "@"
$cD.Value = "115.47"
$cD.Style = "Normal"
$ws.Cells.Item(28, 5).Value = "  +0.91%  "
$cD = $ws.Cells.Item(29, 4)
$cD.NumberFormat = "@"
$cD.Value = "4.931"
$cD.Style = "Normal"
$ws.Cells.Item(29, 5).Value = "  +1.01%  "
$cD = $ws.Cells.Item(30, 4)
$cD.NumberFormat = "@"
$cD.Value = "0.08886"
$cD.Style = "Normal"
$ws.Cells.Item(30, 5).Value = "  +0.27%  "
$cD = $ws.Cells.Item(31, 4)
$cD.NumberFormat = "@"
$cD.Value = "3.312"
$cD.Style = "Normal"
$ws.Cells.Item(31, 5).Value = "  +2.93%  "
$ws.Cells.Item(32, 5).Value = "  +0.40%  "
$cD = $ws.Cells.Item(33, 4)
$cD.NumberFormat = "@"
$cD.Value = "4.577"
$cD.Style = "Normal"
$ws.Cells.Item(33, 5).Value = "  +2.17%  "
$cD = $ws.Cells.Item(34, 4)
$cD.NumberFormat = "@"
$cD.Value = "0.7366"
$cD.Style = "Normal"
$ws.Cells.Item(34, 5).Value = "  -1.37%  "
$cD = $ws.Cells.Item(35, 4)
$cD.NumberFormat = "@"
$cD.Value = "2.689"
$cD.Style = "Normal"
$ws.Cells.Item(35, 5).Value = "  -3.40%  "
$cD = $ws.Cells.Item(36, 4)
$cD.NumberFormat = "@"
$cD.Value = "1.122"
$cD.Style = "Normal"
$ws.Cells.Item(36, 5).Value = "  +3.12%  "
$cD = $ws.Cells.Item(37, 4)
$cD.NumberFormat = "@"
$cD.Value = "0.02007"
$cD.Style = "Normal"
$ws.Cells.Item(37, 5).Value = "  +3.46%  "
$cD = $ws.Cells.Item(38, 4)
$cD.NumberFormat = "@"
$cD.Value = "2.995"
$cD.Style = "Normal"
$ws.Cells.Item(38, 5).Value = "  +0.97%  "
$cD = $ws.Cells.Item(39, 4)
$cD.NumberFormat = "@"
$cD.Value = "0.05252"
$cD.Style = "Normal"
$ws.Cells.Item(39, 5).Value = "  +0.92%  "
$cD = $ws.Cells.Item(40, 4)
$cD.NumberFormat = "@"
$cD.Value = "0.5326"
$cD.Style = "Normal"
$cD = $ws.Cells.Item(41, 4)
$cD.NumberFormat = "@"
$cD.Value = "7.064"
$cD.Style = "Normal"
$ws.Cells.Item(41, 5).Value = "  +1.44%  "
$ws.Cells.Item(42, 5).Value = "  +0.60%  "
$cD = $ws.Cells.Item(43, 4)
$cD.NumberFormat = "@"
$cD.Value = "8.337"
$cD.Style = "Normal"
$ws.Cells.Item(43, 5).Value = "  +2.18%  "
$cD = $ws.Cells.Item(44, 4)
$cD.NumberFormat = "@"
$cD.Value = "10.57"
$cD.Style = "Normal"
$ws.Cells.Item(44, 5).Value = "  +0.97%  "
$cD = $ws.Cells.Item(45, 4)
$cD.NumberFormat = "@"
$cD.Value = "0.4770"
$cD.Style = "Normal"
$ws.Cells.Item(45, 5).Value = "  +1.42%  "
$ws.Cells.Item(46, 5).Value = "  +1.26%  "
$cD = $ws.Cells.Item(47, 4)
$cD.NumberFormat = "@"
$cD.Value = "102.44"
$cD.Style = "Normal"
$ws.Cells.Item(47, 5).Value = "  +1.72%  "
$cD = $ws.Cells.Item(48, 4)
$cD.NumberFormat = "@"
$cD.Value = "1.629"
$cD.Style = "Normal"
$ws.Cells.Item(48, 5).Value = "  +1.83%  "
$cD = $ws.Cells.Item(49, 4)
$cD.NumberFormat = "@"
$cD.Value = "66.41"
$cD.Style = "Normal"
$ws.Cells.Item(49, 5).Value = "  +1.25%  "
$cD = $ws.Cells.Item(51, 4)
$cD.NumberFormat = "@"
$cD.Value = "0.8969"
$cD.Style = "Normal"
$ws.Cells.Item(51, 5).Value = "  -0.05%  "
